$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H-column subcategory labels to pluralized / revised wording.
$ws.Range("H3").Value = "drawing(s)"
$ws.Range("H15").Value = "line graph(s)"
$ws.Range("H22").Value = "line graph(s)"
$ws.Range("H28").Value = "line graph(s)"
$ws.Range("H29").Value = "line graph(s)"
$ws.Range("H39").Value = "photo(s)"
$ws.Range("H40").Value = "photo(s)"
$ws.Range("H41").Value = "drawing(s)"
$ws.Range("H42").Value = "photo(s)"
$ws.Range("H43").Value = "photo(s)"
$ws.Range("H47").Value = "data display"
$ws.Range("H48").Value = "data collection, data analysis, data gathering diagram"

# Remove the now-unused "is_viewed" column (I) entirely, including its
# header, shrinking the sheet's used range back down to column H.
$ws.Range("I1:I55").Delete()
